$pairs = @(
    @('ahi3', 'Apnea / Hypopnea events with >= 3% percent desaturation per hour of sleep'),
    @('ahi_a0h3', 'AHI 3% -- all apneas and hypopneas with >=3% oxygen desaturation'),
    @('ahi_a0h3a', 'AHI 3% -- all apneas and hypopneas with >=3% oxygen desaturation or arousal'),
    @('ahi_a0h4', 'AHI 4% -- all apneas and hypopneas with >=4% oxygen desaturation'),
    @('ahi_a0h4a', 'AHI 4% -- all apneas and hypopneas with >=4% oxygen desaturation or arousal'),
    @('ahi_c0h3', 'Central AHI 3% -- central apneas and hypopneas with >=3% oxygen desaturation'),
    @('ahi_c0h3a', 'Central AHI 3% -- central apneas and hypopneas with >=3% oxygen desaturation or arousal'),
    @('ahi_c0h4', 'Central AHI 4% -- central apneas and hypopneas with >=4% oxygen desaturation'),
    @('ahi_c0h4a', 'Central AHI 4% -- central apneas and hypopneas with >=4% oxygen desaturation or arousal'),
    @('ahi_o0h3', 'Obstructive AHI 3% -- obstructive apneas and hypopneas with >=3% oxygen desaturation'),
    @('ahi_o0h3a', 'Obstructive AHI 3% -- obstructive apneas and hypopneas with >=3% oxygen desaturation or arousal'),
    @('ahi_o0h4', 'Obstructive AHI 4% -- obstructive apneas and hypopneas with >=4% oxygen desaturation'),
    @('ahi_o0h4a', 'Obstructive AHI 4% -- obstructive apneas and hypopneas with >=4% oxygen desaturation or arousal'),
)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Block 1: row 40 originally held the "icsdadultosa_psg5events" / ahi summary row.
# Update it in place to the first (ahi3) pair, then insert 12 new rows below it
# for the remaining pairs, all tagged with the same id (A column).

$ws.Range("B40").Value = $pairs[0][0]
$ws.Range("C40").Value = $pairs[0][1]

for ($i = 1; $i -lt $pairs.Length; $i++) {
    $targetRow = 40 + $i
    $ws.Rows.Item($targetRow).Insert()
    $ws.Cells.Item($targetRow, 1).Value = "icsdadultosa_psg5events"
    $ws.Cells.Item($targetRow, 2).Value = $pairs[$i][0]
    $ws.Cells.Item($targetRow, 3).Value = $pairs[$i][1]
    $ws.Cells.Item($targetRow, 2).Style = "Normal"
    $ws.Cells.Item($targetRow, 3).Style = "Normal"
}

# --- Block 2: the row that used to be 41 ("icsdadultosa_psg15events") has now been
# pushed down to row 53 by the 12 inserted rows above. Repeat the same expansion.

$secondStart = 40 + $pairs.Length
$ws.Range("B$secondStart").Value = $pairs[0][0]
$ws.Range("C$secondStart").Value = $pairs[0][1]

for ($i = 1; $i -lt $pairs.Length; $i++) {
    $targetRow = $secondStart + $i
    $ws.Rows.Item($targetRow).Insert()
    $ws.Cells.Item($targetRow, 1).Value = "icsdadultosa_psg15events"
    $ws.Cells.Item($targetRow, 2).Value = $pairs[$i][0]
    $ws.Cells.Item($targetRow, 3).Value = $pairs[$i][1]
    $ws.Cells.Item($targetRow, 2).Style = "Normal"
    $ws.Cells.Item($targetRow, 3).Style = "Normal"
}

$lastRow = $secondStart + $pairs.Length - 1
$selected = $ws.Range("A" + ($secondStart + 1) + ":A" + $lastRow).Select()
